# 1st commit by VG
# Update the test-data email addresses on LoginSheet and RegSheet, and
# move the active-cell selection on both sheets (cosmetic, matches what
# Excel records when a user simply clicks around before saving).

$wb = $excel.ActiveWorkbook

# --- RegSheet ---
$regSheet = $wb.Worksheets.Item("RegSheet")
$regSheet.Activate()
$regSheet.Range("A2").Value = "abcx26@gmail.com"
$regSheet.Range("A3").Value = "abcx22@gmail.com"
$regSheet.Range("A4").Value = "abcx23@gmail.com"
$regSheet.Range("A5").Value = "abcx24@gmail.com"
$regSheet.Range("C14").Select()

# --- LoginSheet ---
$loginSheet = $wb.Worksheets.Item("LoginSheet")
$loginSheet.Activate()
$loginSheet.Range("A2").Value = "abcx26@gmail.com"
$loginSheet.Range("A3").Value = "abcx22@gmail.com"
$loginSheet.Range("A4").Value = "abcx23@gmail.com"
$loginSheet.Range("A5").Value = "abcx24@gmail.com"
$loginSheet.Range("F15").Select()
